# Simulated Wild Card round and logged it
# Extends the running play-by-play logs (YDS, ST) and updates the
# cumulative game-by-game totals (OFF, DEF, ST, TURNS, PEN) for the new game.

$wb = $excel.ActiveWorkbook

$ydsSheet = $wb.Worksheets.Item("YDS")
$ydsSheet.Range("B2").Value = '35 3 6 1 2 2 11 4 5 2 0 3 4 4 1 4 5 1 -1 10 -4 1 3 3 1 8 1 3 5 2 4 1 -1 9 2 5 2 7 7 2 4 2 26 9 16 1 -1 1 5 -2 12 3 -1 -5 7 -2 4 3 6 1 2 0 1 17 17 2 2 2 0 -1 7 12 -4 0 0 0 -1 1 4 4 12 0 7 0 1 7 4 -1 4 0 -3 3 -3 4 7 5 15 4 0 8 3 1 24 9 3 2 5 1 21 4 4 5 3 1 1 1 11 -1 6 9 1 3 4 1 21 6 6 3 1 0 0 6 5 2 12 5 9 2 32 6 1 3 10 10 0 13 3 0 -2 15 5 2 4 9 6 3 -2 1 6 -1 -2 0 1 3 -1 3 4 1 8 5 6 5 6 3 7 4 7 5 2 6 10 4 6 5 4 -2 6 3 2 14 -1 6 3 2 0 1 0 1 1 2 2 2 0 16 12 1 5 2 -1 6 2 13 3 3 2 8 4 11 12 2 -2 3 6 5 1 4 11 2 5 1 7 -6 1 2 9 0 5 -1 15 9 13 4 0 5 5 4 18 2 3 1 5 1 -1 16 11 5 1 5 22 17 2 5 10 -1 2 3 0 1 1 1 21 3 2 1 2 14 17 2 1 5 13 -1 7 -1 5 12 2 10 0 4 1 6 2 7 0 3 1 9 0 2 2 2 2 2 1 6 5 7 4 2 -1 4 6 10 19 -1 7 14 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2 4 12 1 4 -2 4 -3 -2 13 9 3 2 5 3 6 11 0 18 2 0 9 5 9 0 0 2 1 16 1 7 31 2 0 2 2 1 3 1 9 6 8 11 3 4 1 4 15 7 2 6 9 2 9 7 0 7 5 4 8 2 4 4 7 10 8 -1 4 1 13 5 3 7 4 6 4 8 3 6 7 -1 8 4 15 -1 3 5 1 7 12 9 2 4 6 19 0 4 4 6 5 -1 2 4 1 2 1 16 5 2 2 3 2 3 14 14 0 0 1 8 2 7 1 5 1'
$ydsSheet.Range("C2").Value = '4 7 4 9 8 5 22 6 7 12 8 25 9 7 21 26 9 6 17 16 4 1 6 10 3 8 3 7 11 7 5 5 2 7 -1 2 28 19 3 24 -1 10 10 3 11 8 9 32 3 -1 1 -1 14 7 2 9 6 17 31 13 7 3 10 8 27 7 9 4 4 4 4 11 5 22 9 10 4 5 -5 10 7 4 5 7 17 1 16 15 5 16 7 3 11 4 12 15 9 4 4 21 13 5 10 12 15 1 12 3 21 30 0 9 7 7 15 11 6 9 17 23 4 5 7 16 9 3 3 21 20 1 -3 4 10 13 24 3 10 9 14 22 20 27 0 8 7 3 5 11 3 11 75 14 28 25 3 24 5 4 19 3 15 28 9 11 4 8 4 1 9 20 4 26 3 22 8 46 22 29 28 14 44 2 5 11 33 7 10 4 15 7 14 2 7 17 15 9 2 7 8 12 -1 4 9 41 3 28 3 7 18 2 12 7 11 20 3 12 26 23 4 8 5 1 10 19 17 0 15 3 26 48 11 5 5 17 6 9 9 19 12 6 9 5 18 6 10 10 6 6 15 3 19 2 10 11 4 22 20 16 13 12 38 5 19 5 3 21 10 7 20 9 11 5 41 4 7 7 12 7 4 10 9 4 3 5 15 25 7 9 10 14 10 8 11 13 12 12 9 20 4 12 43 7 9 16 9 13 5 14 6 10 9 3 6 9 13 21 8 7 12 6 14 9 17 8 6 8 8 4 12 20 20 17 8 8 8 7 7 13 46 8 9 2 5 35 7 5 13 23 17 6 2 28 6 32 5 1 18 39 11 8 9 30 4 12 10 2 5 8 7 10 7 1 43 4 7 15 3 6 2 18'
$ydsSheet.Range("B3").Value = '2 15 4 -1 8 3 -1 -1 2 2 1 15 8 4 -3 5 2 5 4 2 8 13 3 6 1 3 7 3 3 -1 17 5 1 4 4 6 10 -3 4 -2 10 17 3 0 0 4 4 14 11 2 1 0 5 0 7 5 3 6 5 3 3 3 3 8 5 2 7 9 5 7 5 -1 1 -2 1 -3 3 11 3 14 6 3 2 5 5 4 21 -1 11 3 4 3 3 6 3 4 6 8 0 1 7 7 3 8 6 1 4 6 4 1 -2 3 2 8 1 1 2 11 1 2 5 3 0 2 1 3 1 2 1 3 7 1 13 3 4 3 2 2 5 0 2 -1 9 4 1 12 3 4 -1 0 0 0 1 3 9 0 2 4 4 21 8 8 10 2 7 1 3 5 5 0 5 1 3 6 3 1 6 -3 19 4 3 6 3 5 3 10 75 3 3 1 9 3 2 28 4 2 2 3 0 2 4 2 0 4 4 4 -2 3 10 2 1 3 9 6 0 3 0 4 10 1 4 4 5 10 24 19 0 3 1 2 5 2 -3 4 5 3 9 1 3 3 5 -2 1 3 -3 3 0 12 2 1 0 2 5 5 0 -1 -1 4 5 2 11 2 5 2 1 4 2 17 3 13 3 68 6 30 5 6 3 7 14 4 2 8 4 3 1 2 1 5 12 9 2 3 1 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1 37 11 4 6 4 5 3 4 2 3 1 11 7 2 3 1 0 9 7 1 6 3 7 -2 2 12 6 3 0 0 -4 1 2 67 3 3 5 1 5 3 25 2 6 2 -1 8 3 7 3 2 11 1 2 2 7 2 0 0 5 4 15 5 2 5 1 5 1 4 2 4 3 4 9 5 4 7 4 2 5 4 3 0 7 5 2 1 10 6 5 3 1 3 3 11 9 1 27 7 3 5 1 3 1 4 3 2 2 23 4 -1 5 2 -2 4 11 3 2 0 4 26 15 1 7 6 5 9 3 8 2 4 9 1 3 9 16 2 15 10 1 3 7 6 4 0'
$ydsSheet.Range("C3").Value = '17 5 9 18 23 6 36 4 0 30 12 3 -2 15 3 13 27 7 11 9 8 5 3 17 4 5 3 27 5 8 6 6 12 12 9 12 17 11 11 6 13 4 9 7 10 12 7 8 6 13 13 28 3 6 0 27 8 1 4 28 26 11 11 9 7 23 18 11 8 8 10 7 11 7 67 9 16 13 37 7 15 5 -5 27 27 27 3 14 8 11 14 11 18 17 3 23 1 13 6 13 44 20 11 3 11 4 13 0 18 27 20 1 13 12 33 3 2 4 12 24 14 6 9 7 35 12 13 7 -2 6 15 3 7 22 5 17 21 8 7 21 8 8 12 15 6 7 3 4 6 6 16 12 41 6 1 13 5 5 9 11 11 15 14 7 6 18 11 14 24 13 6 2 6 5 33 13 -3 6 21 19 17 6 9 7 12 6 13 11 11 -1 9 6 9 2 11 5 18 3 9 5 18 6 12 1 9 8 9 9 7 0 10 0 6 5 11 4 19 16 17 4 8 10 2 14 8 1 2 3 24 7 14 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6 7 8 23 0 19 11 8 7 13 7 8 12 9 4 13 3 10 1 12 9 11 18 9 11 11 6 5 28 9 17 15 19 2 7 9 40 2 8 24 11 22 5 4 3 -5 4 16 6 9 28 9 4 5 13 7 9 7 17 4 3 5 9 2 0 15 9 0 7 8 24 2 19 6 11 22 19 45 4 9 4 34 8 19 19 38 1'

$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("B2").Value = 7
$offSheet.Range("C2").Value = 230
$offSheet.Range("D2").Value = 16
$offSheet.Range("F2").Value = 59
$offSheet.Range("G2").Value = 74
$offSheet.Range("I2").Value = 6
$offSheet.Range("J2").Value = 43
$offSheet.Range("N2").Value = 16
$offSheet.Range("O2").Value = 20
$offSheet.Range("P2").Value = 15
$offSheet.Range("C3").Value = 173
$offSheet.Range("D3").Value = 6
$offSheet.Range("E3").Value = 37
$offSheet.Range("F3").Value = 110
$offSheet.Range("G3").Value = 28
$offSheet.Range("H3").Value = 24
$offSheet.Range("I3").Value = 57
$offSheet.Range("J3").Value = 60
$offSheet.Range("L3").Value = 253
$offSheet.Range("M3").Value = 169
$offSheet.Range("Q3").Value = 547

$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("C2").Value = 221
$defSheet.Range("D2").Value = 15
$defSheet.Range("E2").Value = 9
$defSheet.Range("F2").Value = 71
$defSheet.Range("G2").Value = 71
$defSheet.Range("I2").Value = 8
$defSheet.Range("J2").Value = 41
$defSheet.Range("N2").Value = 19
$defSheet.Range("O2").Value = 32
$defSheet.Range("B3").Value = 12
$defSheet.Range("C3").Value = 178
$defSheet.Range("F3").Value = 112
$defSheet.Range("G3").Value = 41
$defSheet.Range("H3").Value = 30
$defSheet.Range("I3").Value = 57
$defSheet.Range("J3").Value = 55
$defSheet.Range("L3").Value = 249
$defSheet.Range("M3").Value = 151
$defSheet.Range("Q3").Value = 520

$stSheet = $wb.Worksheets.Item("ST")
$stSheet.Range("B2").Value = 107
$stSheet.Range("D2").Value = 53
$stSheet.Range("F2").Value = 440
$stSheet.Range("G2").Value = 428
$stSheet.Range("L2").Value = 122
$stSheet.Range("M2").Value = 91
$stSheet.Range("B3").Value = 44
$stSheet.Range("B4").Value = '54 65 61 61 62 56 63 64 51 68 26 64 64 67 67 60 63 61 59 64 67 64 67 65 67 66 66 61 68 56 54 54 60 66 65 64 65 65 64 54 67 62 64 61 67 63 56 60 65 54 63 62 67 61 62 60 61 62 62 50 58 57 42'
$stSheet.Range("B5").Value = '0 38 18 21 23 16 0 13 10 24 0 26 22 24 25 22 27 27 22 23 25 19 25 17 25 29 21 22 19 20 12 22 18 28 19 25 28 21 16 17 27 18 20 22 30 16 27 21 22 20 18 31 25 18 16 9 20 20 18 15 5 15 0'
$stSheet.Range("B6").Value = '17 23 18 0 25 0 18 21 23 24 31 26 26 30 14 25 21 26 37 20 17 14 27 9 24 15 16 21 25 26'
$stSheet.Range("D3").Value = '37 62 32 48 44 35 60 64 52 58 41 68 37 49 44 59 58 58 44 41 42 34 35 46 38 62 55 42 55 66 44 42 60 47 48 36 15 46 45 39 71 46 34 46 49 63 40 46 45 47 47 42 53'
$stSheet.Range("D4").Value = '0 18 0 6 0 0 0 25 0 0 0 0 0 10 0 5 12 11 -1 0 -1 0 0 0 0 12 0 0 4 0 0 0 11 0 0 8 0 0 0 0 0 0 0 4 15 15 3 0 0 0 0 2 52'
$stSheet.Range("D5").Value = '7 0 -1 14 0 0 8 0 18 12 0 0 11 6 0 13 0 23 17 13 21 22 10 0 27 0 18 15 7 0 0 0 0 0 10 0 0 17 0 15 8 0 0 0 0 1 0 0 0 0 0 11 0 9 0'

$turnsSheet = $wb.Worksheets.Item("TURNS")
$turnsSheet.Range("B3").Value = 8

$penSheet = $wb.Worksheets.Item("PEN")
$penSheet.Range("B2").Value = 17
$penSheet.Range("D2").Value = 10
$penSheet.Range("D4").Value = 8

